$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend dates in row 1 from V1 to AT1 (43045 .. 43069)
$startDate = 43045
for ($i = 0; $i -lt 25; $i++) {
    $col = 22 + $i   # V=22
    $ws.Cells.Item(1, $col).Value = $startDate + $i
    $ws.Cells.Item(1, $col).NumberFormat = "d-mmm"
}
Write-Host "dates done"

# Row 7 text change: Final testing -> Prototype testing
$ws.Range("A7").Value = "Prototype testing"

# New rows 8-14
$ws.Range("A8").Value = "Text chat"
$ws.Range("A9").Value = "Accounts"
$ws.Range("A10").Value = "Create database management"
$ws.Range("A11").Value = "Implement accounts"
$ws.Range("A12").Value = "Recorded sessions"
$ws.Range("A13").Value = "Plan storage framework"
$ws.Range("A14").Value = "Implement recording"
Write-Host "rows text done"

# Bold "category" rows (matches existing pattern used for rows 2,3,6,7)
$ws.Range("A8").Font.Bold = $true
$ws.Range("A9").Font.Bold = $true
$ws.Range("A12").Font.Bold = $true
# Sub-task rows get a distinct (visually identical) explicit-font style
$ws.Range("A10").Font.Bold = $false
$ws.Range("A11").Font.Bold = $false
$ws.Range("A13").Font.Bold = $false
$ws.Range("A14").Font.Bold = $false
Write-Host "bold done"

# ---- Gantt task-duration bars (medium box borders) ----
$blue  = 12611584   # RGB(0,112,192)  FF0070C0
$lblue = 15773696   # RGB(0,176,240)  FF00B0F0
$xlNone = -4142

function Draw-Bar {
    param($range, $color)
    $r = $ws.Range($range)
    $r.BorderAround(1, -4138, 1, $color) | Out-Null
}

function Clear-Bottom {
    param($range)
    $ws.Range($range).Borders.Item(9).LineStyle = $xlNone
}

Draw-Bar "V8:AA8"   $blue    # Text chat
Draw-Bar "AB9:AL9"  $blue    # Accounts
Clear-Bottom "AB9:AF9"       # open segment above "Create database management"
Draw-Bar "AB10:AF10" $lblue  # Create database management
Draw-Bar "AG11:AL11" $lblue  # Implement accounts
Draw-Bar "AM12:AT12" $blue   # Recorded sessions
Clear-Bottom "AM12:AO12"     # open segment above "Plan storage framework"
Draw-Bar "AM13:AO13" $lblue  # Plan storage framework
Draw-Bar "AP14:AT14" $lblue  # Implement recording
Write-Host "bars done"
